$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 0.5876173973083496
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 0.5127014517784119
$ws.Cells.Item(4, 4).Value = 0.8669266700744629
$ws.Cells.Item(5, 4).Value = 0.8246411681175232
$ws.Cells.Item(6, 4).Value = 0.06286820024251938
$ws.Cells.Item(7, 4).Value = 0.8776209354400635
$ws.Cells.Item(8, 4).Value = 0.7707201838493347
$ws.Cells.Item(9, 4).Value = 0.3048694133758545
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0.495371550321579
$ws.Cells.Item(11, 4).Value = 0.2303283959627151
$ws.Cells.Item(12, 4).Value = 0.1857190877199173
$ws.Cells.Item(13, 4).Value = 0.8674956560134888
$ws.Cells.Item(14, 4).Value = 0.8751955032348633
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0.4216914474964142
$ws.Cells.Item(16, 4).Value = 0.2624797523021698
$ws.Cells.Item(17, 4).Value = 0.7577654123306274
$ws.Cells.Item(18, 4).Value = 0.6054984331130981
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 0.4093546569347382
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 0.3960031867027283
$ws.Cells.Item(21, 4).Value = 0.1320594996213913
$ws.Cells.Item(22, 4).Value = 0.6615296602249146
$ws.Cells.Item(23, 4).Value = 0.7411167621612549
$ws.Cells.Item(24, 4).Value = 0.4433744549751282
$ws.Cells.Item(25, 4).Value = 0.7918363809585571
$ws.Cells.Item(26, 4).Value = 0.9470779895782471
$ws.Cells.Item(27, 4).Value = 0.04427319392561913
$ws.Cells.Item(28, 4).Value = 0.8608905076980591
$ws.Cells.Item(29, 4).Value = 0.3152673542499542
$ws.Cells.Item(30, 4).Value = 0.7652621269226074
$ws.Cells.Item(31, 4).Value = 0.5419816970825195
$ws.Cells.Item(32, 4).Value = 0.7014244794845581
$ws.Cells.Item(33, 4).Value = 0.722465991973877
$ws.Cells.Item(34, 4).Value = 0.907707154750824
$ws.Cells.Item(35, 4).Value = 0.6762729287147522
$ws.Cells.Item(36, 4).Value = 0.1494694650173187
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 0.453951358795166
$ws.Cells.Item(38, 4).Value = 0.8986181616783142
$ws.Cells.Item(39, 4).Value = 0.9582871198654175
$ws.Cells.Item(40, 4).Value = 0.8870537877082825
$ws.Cells.Item(41, 4).Value = 0.9650564789772034
$ws.Cells.Item(42, 4).Value = 0.5881452560424805
$ws.Cells.Item(43, 4).Value = 0.7852417826652527
$ws.Cells.Item(44, 4).Value = 0.6357388496398926
$ws.Cells.Item(45, 4).Value = 0.4414702951908112
$ws.Cells.Item(46, 4).Value = 0.8643794059753418
$ws.Cells.Item(47, 4).Value = 0.7324144244194031
$ws.Cells.Item(48, 4).Value = 0.7587654590606689
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 0.3502878844738007
$ws.Cells.Item(50, 4).Value = 0.4270217418670654
$ws.Cells.Item(51, 4).Value = 0.7436496019363403
$ws.Cells.Item(52, 4).Value = 0.08126440644264221
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 0.2828573286533356
$ws.Cells.Item(54, 4).Value = 0.8288872241973877
$ws.Cells.Item(55, 4).Value = 0.3590312004089355
$ws.Cells.Item(56, 3).Value = 1
$ws.Cells.Item(56, 4).Value = 0.6459195613861084
$ws.Cells.Item(57, 4).Value = 0.2269430160522461
$ws.Cells.Item(58, 4).Value = 0.8229509592056274
$ws.Cells.Item(59, 4).Value = 0.7895099520683289
$ws.Cells.Item(60, 4).Value = 0.8976601958274841
$ws.Cells.Item(61, 4).Value = 0.6911138892173767
$ws.Cells.Item(62, 4).Value = 0.7324249148368835
$ws.Cells.Item(63, 4).Value = 0.5695649385452271
$ws.Cells.Item(64, 4).Value = 0.3720881342887878
$ws.Cells.Item(65, 4).Value = 0.09626663476228714
$ws.Cells.Item(66, 4).Value = 0.746989369392395
$ws.Cells.Item(67, 3).Value = 0
$ws.Cells.Item(67, 4).Value = 0.1907045394182205
$ws.Cells.Item(68, 4).Value = 0.590224027633667
$ws.Cells.Item(69, 4).Value = 0.6740698218345642
$ws.Cells.Item(70, 4).Value = 0.21464604139328
$ws.Cells.Item(71, 4).Value = 0.2641399502754211
$ws.Cells.Item(72, 4).Value = 0.277616560459137
$ws.Cells.Item(73, 4).Value = 0.8923956155776978
$ws.Cells.Item(74, 3).Value = 1
$ws.Cells.Item(74, 4).Value = 0.6170204281806946
$ws.Cells.Item(75, 4).Value = 0.1127720177173615
$ws.Cells.Item(76, 4).Value = 0.4118671715259552
$ws.Cells.Item(77, 4).Value = 0.6918544769287109
$ws.Cells.Item(78, 4).Value = 0.8272613286972046
$ws.Cells.Item(79, 4).Value = 0.1887997835874557
